$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "Μπανάνες" (banana) rows entirely (old rows 11 then 4;
# delete bottom row first so the earlier row number stays valid).
$ws.Rows(11).Delete()
$ws.Rows(4).Delete()

# After the two deletions the rows are (1-indexed):
#  2 Nivea (unchanged)
#  3 Βερύκοκα/Apricot  -> should become Ροδάκινα/Peach
#  4 Παπαγάλος/Coffee  -> should become Ροδάκινα/Peach (coffee row moves to row 7)
#  5 Ροδάκινα/Peach    -> should become Βερύκοκα/Apricot
#  6 Ροδάκινα/Peach    -> stays Ροδάκινα/Peach
#  7 Ροδάκινα/Peach    -> should become Παπαγάλος/Coffee
#  8 Βερύκοκα/Apricot  -> stays Βερύκοκα/Apricot
#  9 Nutella           -> stays Nutella
# 10 Νεκταρίνια        -> stays Νεκταρίνια
# 11 totals            -> stays totals

$cols = 1,5,6,7,8,9,10

# Move the coffee row (currently row 4) down to row 7.
foreach ($col in $cols) {
    $ws.Cells.Item(4, $col).Copy($ws.Cells.Item(7, $col))
}

# Turn rows 3 and 4 into peach rows (copy from row 6, a stable peach row).
foreach ($col in $cols) {
    $ws.Cells.Item(6, $col).Copy($ws.Cells.Item(3, $col))
    $ws.Cells.Item(6, $col).Copy($ws.Cells.Item(4, $col))
}

# Turn row 5 into an apricot row (copy from row 8, a stable apricot row).
foreach ($col in $cols) {
    $ws.Cells.Item(8, $col).Copy($ws.Cells.Item(5, $col))
}

# --- Update the quantity (K) / turnover (L) figures ---
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 3.84

$ws.Cells.Item(6, 11).Value = 3.535
$ws.Cells.Item(6, 12).Value = 2.66

$ws.Cells.Item(7, 11).Value = 6
$ws.Cells.Item(7, 12).Value = 9.81

$ws.Cells.Item(8, 11).Value = 6.9
$ws.Cells.Item(8, 12).Value = 10.07

$ws.Cells.Item(9, 11).Value = 12
$ws.Cells.Item(9, 12).Value = 32.06

$ws.Cells.Item(10, 11).Value = 40.129
$ws.Cells.Item(10, 12).Value = 30.8

# Totals row
$ws.Cells.Item(11, 11).Value = 71.564
$ws.Cells.Item(11, 12).Value = 89.24

# Shrink the conditional-formatting ranges to match the new data extent.
$ws.Range("I1:I12").FormatConditions.Delete()
$ws.Range("J1:J12").FormatConditions.Delete()

$cf1 = $ws.Range("I1:I10").FormatConditions.AddColorScale(3)
$cf1.ColorScaleCriteria.Item(1).Type = 1
$cf1.ColorScaleCriteria.Item(1).FormatColor.Color = 7039083
$cf1.ColorScaleCriteria.Item(2).Type = 4
$cf1.ColorScaleCriteria.Item(2).Value = 50
$cf1.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167
$cf1.ColorScaleCriteria.Item(3).Type = 2
$cf1.ColorScaleCriteria.Item(3).FormatColor.Color = 8105619

$cf2 = $ws.Range("J1:J10").FormatConditions.AddColorScale(3)
$cf2.ColorScaleCriteria.Item(1).Type = 1
$cf2.ColorScaleCriteria.Item(1).FormatColor.Color = 7039083
$cf2.ColorScaleCriteria.Item(2).Type = 4
$cf2.ColorScaleCriteria.Item(2).Value = 50
$cf2.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167
$cf2.ColorScaleCriteria.Item(3).Type = 2
$cf2.ColorScaleCriteria.Item(3).FormatColor.Color = 8105619
